$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 status update ---
$ws.Range("C6").Value = "Resolved"

# --- New rows 23-28 (content first, so new shared strings are appended
#     in the same order the reference workbook used them) ---
$ws.Range("A23").Value = 42614
$ws.Range("B23").Value = "Remove last column onclick in melody cells"
$ws.Range("C23").Value = "Resolved"

$ws.Range("A24").Value = 42614
$ws.Range("B24").Value = "Right side key is not written"
$ws.Range("C24").Value = "Resolved"

$ws.Range("A25").Value = 42614
$ws.Range("B25").Value = "Left to right modulation is not stopping properly"
$ws.Range("C25").Value = "Pending"

$ws.Range("A26").Value = 42614
$ws.Range("B26").Value = "Add major and minor key indicator"
$ws.Range("C26").Value = "Pending"

$ws.Range("A27").Value = 42614
$ws.Range("B27").Value = "replace lower rows and register"
$ws.Range("C27").Value = "Pending"

$ws.Range("D25").Value = "Add leading tone rule"

$ws.Range("A28").Value = 42619
$ws.Range("B28").Value = "add 2 octaves"
$ws.Range("C28").Value = "Pending"

# --- Row 6 note update (last, so "works on Chrome" lands at the final index) ---
$ws.Range("D6").Value = "works on Chrome"

# --- Apply date styling (copy format from A2) to new A cells ---
$ws.Range("A23:A28").NumberFormat = $ws.Range("A2").NumberFormat

# --- Apply yellow fill (copy format from B4, style s=3) to new B cells ---
$ws.Range("B23:B28").Interior.Color = $ws.Range("B4").Interior.Color

# --- Sheet view changes ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D6").Select()
